$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17, shifting existing rows 17-20 down to 18-21.
# Excel will copy formatting from the row above (row 16), keeping the date
# style (s="2") consistent on column D.
$ws.Rows("17:17").Insert()

# The row that used to be row 16 (now still row 16) gets brand-new weekly
# price data, while the freshly inserted row 17 receives the values that
# previously lived in row 16 (i.e. the old row 16 record moved down one row).

# Row 16 - new weekly entry
$ws.Cells.Item(16, 1).Value = 11
$ws.Cells.Item(16, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(16, 3).Value = "Bíobío"
$ws.Cells.Item(16, 4).Value = 44511
$ws.Cells.Item(16, 5).Value = 8
$ws.Cells.Item(16, 6).Value = "Fruta"
$ws.Cells.Item(16, 7).Value = 100107
$ws.Cells.Item(16, 8).Value = "Otros"
$ws.Cells.Item(16, 9).Value = 100107002
$ws.Cells.Item(16, 10).Value = "Chirimoya"
$ws.Cells.Item(16, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(16, 12).Value = "Primera"
$ws.Cells.Item(16, 13).Value = 80
$ws.Cells.Item(16, 14).Value = 25000
$ws.Cells.Item(16, 15).Value = 26000
$ws.Cells.Item(16, 16).Value = 25375
$ws.Cells.Item(16, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(16, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(16, 19).Value = 2538
$ws.Cells.Item(16, 20).Value = 10

# Row 17 - carries forward the data that used to be in row 16
$ws.Cells.Item(17, 1).Value = 11
$ws.Cells.Item(17, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(17, 3).Value = "Bíobío"
$ws.Cells.Item(17, 4).Value = 44469
$ws.Cells.Item(17, 5).Value = 8
$ws.Cells.Item(17, 6).Value = "Fruta"
$ws.Cells.Item(17, 7).Value = 100107
$ws.Cells.Item(17, 8).Value = "Otros"
$ws.Cells.Item(17, 9).Value = 100107002
$ws.Cells.Item(17, 10).Value = "Chirimoya"
$ws.Cells.Item(17, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(17, 12).Value = "Primera"
$ws.Cells.Item(17, 13).Value = 100
$ws.Cells.Item(17, 14).Value = 28000
$ws.Cells.Item(17, 15).Value = 29000
$ws.Cells.Item(17, 16).Value = 28500
$ws.Cells.Item(17, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(17, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(17, 19).Value = 2850
$ws.Cells.Item(17, 20).Value = 10

# Ensure the date cells keep the datetime number format used throughout
# column D (style index referencing numFmtId 165), matching row 16's style.
$ws.Cells.Item(17, 4).NumberFormat = $ws.Cells.Item(16, 4).NumberFormat
